$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C18").Value = "PB12"
